$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Mobile networks")

# Insert a new row above row 3 (before "Sky Mobile") to make room for "O2",
# keeping the list in alphabetical order: EE, O2, Sky Mobile, SMARTY, ...
$ws.Rows.Item(3).Insert()
$ws.Range("A3").Value = "O2"

# Move selection to match the authored workbook state
$ws.Range("A3").Select()
